$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 106, shifting existing rows 106-115 down to 107-116.
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with the new record (copy of the constant columns
# A-C, E-J plus the new record's own D, K, L, M, N, O, P, Q, R, S, T values).
$ws.Cells.Item(106, 1).Value = 1
$ws.Cells.Item(106, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(106, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(106, 4).Value = 44783
$ws.Cells.Item(106, 5).Value = 15
$ws.Cells.Item(106, 6).Value = "Fruta"
$ws.Cells.Item(106, 7).Value = 100106
$ws.Cells.Item(106, 8).Value = "Oleaginosos"
$ws.Cells.Item(106, 9).Value = 100106002
$ws.Cells.Item(106, 10).Value = "Palta"
$ws.Cells.Item(106, 11).Value = "Fuerte"
$ws.Cells.Item(106, 12).Value = "Tercera"
$ws.Cells.Item(106, 13).Value = 200
$ws.Cells.Item(106, 14).Value = 45000
$ws.Cells.Item(106, 15).Value = 46000
$ws.Cells.Item(106, 16).Value = 45500
$ws.Cells.Item(106, 17).Value = "`$/caja 25 kilos"
$ws.Cells.Item(106, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(106, 19).Value = 1820
$ws.Cells.Item(106, 20).Value = 25
